$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 233.95  # ALC!H55
$ws.Cells.Item(55, 9).Value = 293.16666  # ALC!I55
$ws.Cells.Item(55, 10).Value = 145.125  # ALC!J55
$ws.Cells.Item(55, 11).Value = 293.16666  # ALC!K55
$ws.Cells.Item(55, 12).Value = 145.125  # ALC!L55
$ws.Cells.Item(55, 13).Value = -79.16665999999998  # ALC!M55
$ws.Cells.Item(55, 14).Value = -573.125  # ALC!N55
$ws.Cells.Item(98, 8).Value = 719.65515  # ALC!H98
$ws.Cells.Item(98, 9).Value = 655.8570999999999  # ALC!I98
$ws.Cells.Item(98, 10).Value = 2506  # ALC!J98
$ws.Cells.Item(98, 11).Value = 655.8570999999999  # ALC!K98
$ws.Cells.Item(98, 12).Value = 2506  # ALC!L98
$ws.Cells.Item(98, 13).Value = 842.1429000000001  # ALC!M98
$ws.Cells.Item(98, 14).Value = -5502  # ALC!N98
$ws.Cells.Item(107, 8).Value = 3837.8333  # ALC!H107
$ws.Cells.Item(107, 9).Value = 1563.1666  # ALC!I107
$ws.Cells.Item(107, 11).Value = 1563.1666  # ALC!K107
$ws.Cells.Item(107, 13).Value = 356.8334  # ALC!M107
$ws.Cells.Item(122, 8).Value = 719.65515  # ALC!H122
$ws.Cells.Item(122, 9).Value = 655.8570999999999  # ALC!I122
$ws.Cells.Item(122, 10).Value = 2506  # ALC!J122
$ws.Cells.Item(122, 11).Value = 1967.5713  # ALC!K122
$ws.Cells.Item(122, 12).Value = 7518  # ALC!L122
$ws.Cells.Item(122, 13).Value = 482.4287000000002  # ALC!M122
$ws.Cells.Item(122, 14).Value = -12418  # ALC!N122
$ws.Cells.Item(135, 8).Value = 13338311  # ALC!H135
$ws.Cells.Item(135, 9).Value = 15156595  # ALC!I135
$ws.Cells.Item(135, 10).Value = 4223  # ALC!J135
$ws.Cells.Item(135, 11).Value = 136409355  # ALC!K135
$ws.Cells.Item(135, 12).Value = 38007  # ALC!L135
$ws.Cells.Item(135, 13).Value = -136406820  # ALC!M135
$ws.Cells.Item(135, 14).Value = -43077  # ALC!N135
$ws.Cells.Item(137, 8).Value = 941348.5  # ALC!H137
$ws.Cells.Item(137, 9).Value = 43669.832  # ALC!I137
$ws.Cells.Item(137, 11).Value = 131009.496  # ALC!K137
$ws.Cells.Item(137, 13).Value = -128459.496  # ALC!M137
$ws.Cells.Item(138, 8).Value = 4382.9136  # ALC!H138
$ws.Cells.Item(138, 10).Value = 5222.6445  # ALC!J138
$ws.Cells.Item(138, 12).Value = 15667.9335  # ALC!L138
$ws.Cells.Item(138, 14).Value = -25947.9335  # ALC!N138
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 45456596  # ARM!H2
$ws.Cells.Item(2, 9).Value = 58825716  # ARM!I2
$ws.Cells.Item(2, 10).Value = 1592.6  # ARM!J2
$ws.Cells.Item(2, 11).Value = 58825716  # ARM!K2
$ws.Cells.Item(2, 12).Value = 1592.6  # ARM!L2
$ws.Cells.Item(2, 13).Value = -58825603  # ARM!M2
$ws.Cells.Item(2, 14).Value = -1818.6  # ARM!N2
$ws.Cells.Item(32, 8).Value = 26871.812  # ARM!H32
$ws.Cells.Item(32, 9).Value = 15919.551  # ARM!I32
$ws.Cells.Item(32, 10).Value = 95572.37  # ARM!J32
$ws.Cells.Item(32, 11).Value = 15919.551  # ARM!K32
$ws.Cells.Item(32, 12).Value = 95572.37  # ARM!L32
$ws.Cells.Item(32, 13).Value = -15632.551  # ARM!M32
$ws.Cells.Item(32, 14).Value = -96146.37  # ARM!N32
$ws.Cells.Item(61, 8).Value = 3682.8667  # ARM!H61
$ws.Cells.Item(61, 9).Value = 3397.96  # ARM!I61
$ws.Cells.Item(61, 11).Value = 3397.96  # ARM!K61
$ws.Cells.Item(61, 13).Value = -3185.96  # ARM!M61
$ws.Cells.Item(116, 8).Value = 45456596  # ARM!H116
$ws.Cells.Item(116, 9).Value = 58825716  # ARM!I116
$ws.Cells.Item(116, 10).Value = 1592.6  # ARM!J116
$ws.Cells.Item(116, 11).Value = 58825716  # ARM!K116
$ws.Cells.Item(116, 12).Value = 1592.6  # ARM!L116
$ws.Cells.Item(116, 13).Value = -58823422  # ARM!M116
$ws.Cells.Item(116, 14).Value = -6180.6  # ARM!N116
$ws.Cells.Item(122, 8).Value = 1375  # ARM!H122
$ws.Cells.Item(122, 9).Value = 900  # ARM!I122
$ws.Cells.Item(122, 11).Value = 2700  # ARM!K122
$ws.Cells.Item(122, 13).Value = -250  # ARM!M122
$ws.Cells.Item(136, 8).Value = 3682.8667  # ARM!H136
$ws.Cells.Item(136, 9).Value = 3397.96  # ARM!I136
$ws.Cells.Item(136, 11).Value = 10193.88  # ARM!K136
$ws.Cells.Item(136, 13).Value = -7643.880000000001  # ARM!M136
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 45456596  # BSM!H3
$ws.Cells.Item(3, 9).Value = 58825716  # BSM!I3
$ws.Cells.Item(3, 10).Value = 1592.6  # BSM!J3
$ws.Cells.Item(3, 11).Value = 58825716  # BSM!K3
$ws.Cells.Item(3, 12).Value = 1592.6  # BSM!L3
$ws.Cells.Item(3, 13).Value = -58825602  # BSM!M3
$ws.Cells.Item(3, 14).Value = -1820.6  # BSM!N3
$ws.Cells.Item(33, 8).Value = 18000  # BSM!H33
$ws.Cells.Item(33, 9).Value = 18000  # BSM!I33
$ws.Cells.Item(33, 10).Value = 0  # BSM!J33
$ws.Cells.Item(33, 11).Value = 18000  # BSM!K33
$ws.Cells.Item(33, 12).Value = 0  # BSM!L33
$ws.Cells.Item(33, 13).Value = -17664  # BSM!M33
$ws.Cells.Item(33, 14).ClearContents()  # BSM!N33
$ws.Cells.Item(134, 8).Value = 7968.2  # BSM!H134
$ws.Cells.Item(134, 9).Value = 2865.6875  # BSM!I134
$ws.Cells.Item(134, 10).Value = 28378.25  # BSM!J134
$ws.Cells.Item(134, 11).Value = 8597.0625  # BSM!K134
$ws.Cells.Item(134, 12).Value = 85134.75  # BSM!L134
$ws.Cells.Item(134, 13).Value = -6062.0625  # BSM!M134
$ws.Cells.Item(134, 14).Value = -90204.75  # BSM!N134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 11599.6  # CRP!H29
$ws.Cells.Item(29, 10).Value = 11599.6  # CRP!J29
$ws.Cells.Item(29, 12).Value = 11599.6  # CRP!L29
$ws.Cells.Item(29, 14).Value = -12185.6  # CRP!N29
$ws.Cells.Item(31, 8).Value = 9659083  # CRP!H31
$ws.Cells.Item(31, 9).Value = 4221887  # CRP!I31
$ws.Cells.Item(31, 11).Value = 4221887  # CRP!K31
$ws.Cells.Item(31, 13).Value = -4221592  # CRP!M31
$ws.Cells.Item(34, 8).Value = 9659083  # CRP!H34
$ws.Cells.Item(34, 9).Value = 4221887  # CRP!I34
$ws.Cells.Item(34, 11).Value = 4221887  # CRP!K34
$ws.Cells.Item(34, 13).Value = -4221685  # CRP!M34
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 13599.8  # CUL!H70
$ws.Cells.Item(70, 9).Value = 9333  # CUL!I70
$ws.Cells.Item(70, 11).Value = 27999  # CUL!K70
$ws.Cells.Item(70, 13).Value = -27684  # CUL!M70
$ws.Cells.Item(73, 8).Value = 13599.8  # CUL!H73
$ws.Cells.Item(73, 9).Value = 9333  # CUL!I73
$ws.Cells.Item(73, 11).Value = 27999  # CUL!K73
$ws.Cells.Item(73, 13).Value = -26907  # CUL!M73
$ws.Cells.Item(75, 8).Value = 4957.8  # CUL!H75
$ws.Cells.Item(75, 9).Value = 1197.25  # CUL!I75
$ws.Cells.Item(75, 10).Value = 20000  # CUL!J75
$ws.Cells.Item(75, 11).Value = 3591.75  # CUL!K75
$ws.Cells.Item(75, 12).Value = 60000  # CUL!L75
$ws.Cells.Item(75, 13).Value = -2593.75  # CUL!M75
$ws.Cells.Item(75, 14).Value = -61996  # CUL!N75
$ws.Cells.Item(78, 8).Value = 4957.8  # CUL!H78
$ws.Cells.Item(78, 9).Value = 1197.25  # CUL!I78
$ws.Cells.Item(78, 10).Value = 20000  # CUL!J78
$ws.Cells.Item(78, 11).Value = 10775.25  # CUL!K78
$ws.Cells.Item(78, 12).Value = 180000  # CUL!L78
$ws.Cells.Item(78, 13).Value = -5783.25  # CUL!M78
$ws.Cells.Item(78, 14).Value = -189984  # CUL!N78
$ws.Cells.Item(86, 8).Value = 546  # CUL!H86
$ws.Cells.Item(86, 10).Value = 750  # CUL!J86
$ws.Cells.Item(86, 12).Value = 2250  # CUL!L86
$ws.Cells.Item(86, 14).Value = -4622  # CUL!N86
$ws.Cells.Item(89, 8).Value = 546  # CUL!H89
$ws.Cells.Item(89, 10).Value = 750  # CUL!J89
$ws.Cells.Item(89, 12).Value = 6750  # CUL!L89
$ws.Cells.Item(89, 14).Value = -18606  # CUL!N89
$ws.Cells.Item(136, 8).Value = 4840.7  # CUL!H136
$ws.Cells.Item(136, 9).Value = 486.85715  # CUL!I136
$ws.Cells.Item(136, 10).Value = 14999.667  # CUL!J136
$ws.Cells.Item(136, 11).Value = 1460.57145  # CUL!K136
$ws.Cells.Item(136, 12).Value = 44999.001  # CUL!L136
$ws.Cells.Item(136, 13).Value = 3639.42855  # CUL!M136
$ws.Cells.Item(136, 14).Value = -55199.001  # CUL!N136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2048  # GSM!H80
$ws.Cells.Item(80, 9).Value = 2052.8  # GSM!I80
$ws.Cells.Item(80, 10).Value = 2036  # GSM!J80
$ws.Cells.Item(80, 11).Value = 2052.8  # GSM!K80
$ws.Cells.Item(80, 12).Value = 2036  # GSM!L80
$ws.Cells.Item(80, 13).Value = -1054.8  # GSM!M80
$ws.Cells.Item(80, 14).Value = -4032  # GSM!N80
$ws.Cells.Item(83, 8).Value = 2048  # GSM!H83
$ws.Cells.Item(83, 9).Value = 2052.8  # GSM!I83
$ws.Cells.Item(83, 10).Value = 2036  # GSM!J83
$ws.Cells.Item(83, 11).Value = 10264  # GSM!K83
$ws.Cells.Item(83, 12).Value = 10180  # GSM!L83
$ws.Cells.Item(83, 13).Value = -5272  # GSM!M83
$ws.Cells.Item(83, 14).Value = -20164  # GSM!N83
$ws.Cells.Item(122, 8).Value = 9446.697  # GSM!H122
$ws.Cells.Item(122, 10).Value = 3656.8462  # GSM!J122
$ws.Cells.Item(122, 12).Value = 10970.5386  # GSM!L122
$ws.Cells.Item(122, 14).Value = -15870.5386  # GSM!N122
$ws.Cells.Item(141, 8).Value = 94950  # GSM!H141
$ws.Cells.Item(141, 10).Value = 94950  # GSM!J141
$ws.Cells.Item(141, 12).Value = 94950  # GSM!L141
$ws.Cells.Item(141, 14).Value = -105310  # GSM!N141
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2439.875  # LTW!H22
$ws.Cells.Item(22, 9).Value = 2095  # LTW!I22
$ws.Cells.Item(22, 10).Value = 2554.8333  # LTW!J22
$ws.Cells.Item(22, 11).Value = 2095  # LTW!K22
$ws.Cells.Item(22, 12).Value = 2554.8333  # LTW!L22
$ws.Cells.Item(22, 13).Value = -1800  # LTW!M22
$ws.Cells.Item(22, 14).Value = -3144.8333  # LTW!N22
$ws.Cells.Item(27, 8).Value = 2439.875  # LTW!H27
$ws.Cells.Item(27, 9).Value = 2095  # LTW!I27
$ws.Cells.Item(27, 10).Value = 2554.8333  # LTW!J27
$ws.Cells.Item(27, 11).Value = 2095  # LTW!K27
$ws.Cells.Item(27, 12).Value = 2554.8333  # LTW!L27
$ws.Cells.Item(27, 13).Value = -1988  # LTW!M27
$ws.Cells.Item(27, 14).Value = -2768.8333  # LTW!N27
$ws.Cells.Item(46, 8).Value = 1015.8333  # LTW!H46
$ws.Cells.Item(46, 10).Value = 0  # LTW!J46
$ws.Cells.Item(46, 12).Value = 0  # LTW!L46
$ws.Cells.Item(46, 14).ClearContents()  # LTW!N46
$ws.Cells.Item(122, 8).Value = 4941.75  # LTW!H122
$ws.Cells.Item(122, 9).Value = 4383.6665  # LTW!I122
$ws.Cells.Item(122, 11).Value = 13150.9995  # LTW!K122
$ws.Cells.Item(122, 13).Value = -10700.9995  # LTW!M122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1850.3  # WVR!H136
$ws.Cells.Item(136, 9).Value = 1679.5385  # WVR!I136
$ws.Cells.Item(136, 10).Value = 2167.4285  # WVR!J136
$ws.Cells.Item(136, 11).Value = 5038.6155  # WVR!K136
$ws.Cells.Item(136, 12).Value = 6502.2855  # WVR!L136
$ws.Cells.Item(136, 13).Value = -2488.6155  # WVR!M136
$ws.Cells.Item(136, 14).Value = -11602.2855  # WVR!N136

Write-Host "Applied 35 row updates across 8 sheets"
